# Revert "Merge branch 'alpha-michael' into 'master'"
# Restore the "Risk status" / "New Issue needed" columns for the
# "Ethics Application" and "User testing" rows back to Open/No.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 10 ("Ethics Application"): Risk status "Closed" -> "Open",
# New Issue needed "Yes" -> "No"
$t.Cell(10, 10).Range.Text = "Open"
$t.Cell(10, 11).Range.Text = "No"

# Row 11 ("User testing"): Risk status "Closed" -> "Open " (trailing space),
# New Issue needed "Yes" -> "No"
$t.Cell(11, 10).Range.Text = "Open "
$t.Cell(11, 11).Range.Text = "No"
